# Add 2022-Q4 data: a new sheet with fund-level detail, plus a summary row
# on the "总计" (totals) sheet. The new "2022-Q4" sheet is inserted right
# after "总计" and before the existing "2022-Q3" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" sheet (this keeps all the original formatting/styles:
#    header row bold+border, index column border, etc.) and then
#    overwrite its contents with the 2022-Q4 numbers.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Q3 sheet had 8 data rows (rows 2-9); Q4 only needs 5 data rows (rows 2-6)
# so drop the extra three rows at the bottom.
$q4.Rows.Item(7).Resize(3).Delete()

# Force text columns to stay text (avoid "1.60" -> 1.6, "000179" -> 179, etc.)
$q4.Range("B2:G6").NumberFormat = "@"

# Row 2: 000179
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "000179"
$q4.Cells.Item(2,3).Value = "广发美国房地产指数（QDII）人民币A"
$q4.Cells.Item(2,4).Value = "1.60"
$q4.Cells.Item(2,5).Value = "92.49"
$q4.Cells.Item(2,6).Value = "3.02"
$q4.Cells.Item(2,7).Value = "0.0483"
$q4.Cells.Item(2,8).Value = 6

# Row 3: 000180
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "000180"
$q4.Cells.Item(3,3).Value = "广发美国房地产指数（QDII）美元A"
$q4.Cells.Item(3,4).Value = "1.60"
$q4.Cells.Item(3,5).Value = "92.49"
$q4.Cells.Item(3,6).Value = "3.02"
$q4.Cells.Item(3,7).Value = "0.0483"
$q4.Cells.Item(3,8).Value = 6

# Row 4: 070031
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "070031"
$q4.Cells.Item(4,3).Value = "嘉实全球房地产（QDII）"
$q4.Cells.Item(4,4).Value = "0.39"
$q4.Cells.Item(4,5).Value = "94.72"
$q4.Cells.Item(4,6).Value = "2.43"
$q4.Cells.Item(4,7).Value = "0.0095"
$q4.Cells.Item(4,8).Value = 9

# Row 5: 016278
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "016278"
$q4.Cells.Item(5,3).Value = "广发美国房地产指数（QDII）人民币C"
$q4.Cells.Item(5,4).Value = "0.01"
$q4.Cells.Item(5,5).Value = "92.49"
$q4.Cells.Item(5,6).Value = "3.02"
$q4.Cells.Item(5,7).Value = "0.0003"
$q4.Cells.Item(5,8).Value = 6

# Row 6: 016279
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "016279"
$q4.Cells.Item(6,3).Value = "广发美国房地产指数（QDII）美元C"
$q4.Cells.Item(6,4).Value = "0.01"
$q4.Cells.Item(6,5).Value = "92.49"
$q4.Cells.Item(6,6).Value = "3.02"
$q4.Cells.Item(6,7).Value = "0.0003"
$q4.Cells.Item(6,8).Value = 6

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a 2022-Q4 summary row above
#    the existing 2022-Q3 / 2022-Q2 rows (so it reads Q4, Q3, Q2).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Duplicate the border/format of the existing index column (A3, which
# already carries the bordered data-row style) down onto the brand new
# row 4 before we populate it, so its formatting matches rows 2 and 3.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(4,1).PasteSpecial(-4122)

# Shift the old rows down one: old row 3 (2022-Q2) -> row 4,
# old row 2 (2022-Q3) -> row 3, then write the new 2022-Q4 row into row 2.
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 5
$total.Cells.Item(4,4).Value = 0.43

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 8
$total.Cells.Item(3,4).Value = 0.47

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 5
$total.Cells.Item(2,4).Value = 0.11

# Keep the workbook's active tab on "总计", matching the file's original
# bookViews state (activeTab="0").
$total.Activate()
